$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "n"
$ws.Range("B1").Value = "x_m"
$ws.Range("C1").Value = "x_i"
$ws.Range("D1").Value = "x_s"
$ws.Range("E1").Value = "f_m"
$ws.Range("F1").Value = "f_i"
$ws.Range("G1").Value = "f_s"
$ws.Range("H1").Value = "E"

# Copy the existing header style from A1 onto the newly used header cells
# so they share the same bold/centered/bordered style.
$ws.Range("A1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data row 2 ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = -20
$ws.Range("C2").Value = -20
$ws.Range("D2").Value = 40
$ws.Range("E2").Value = 389.999999999999
$ws.Range("F2").Value = 389.999999999999
$ws.Range("G2").Value = [double]"-1.20892581961463e+24"
$ws.Range("H2").Value = 1.000005

# --- Data row 3 ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = -20
$ws.Range("C3").Value = -20
$ws.Range("D3").Value = 40
$ws.Range("E3").Value = 389.999999999999
$ws.Range("F3").Value = 389.999999999999
$ws.Range("G3").Value = [double]"-1.20892581961463e+24"
$ws.Range("H3").Value = 0
